# Realestate Update resale numbers 2024-01-25 21:33
# Append a new data row (row 93) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

# Columns A-D hold text that looks numeric/date-like ("2024-01-25", "03", ...).
# Force them to be stored as text (not auto-converted to a date serial /
# number) by temporarily applying a text number format, then clear the
# format again afterwards so no stray style is left behind on the cell -
# matching the rest of the sheet, where data rows carry no cell style.
$cols = 1,2,3,4
foreach ($col in $cols) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
}

$ws.Cells.Item($row, 1).Value = "2024-01-25"
$ws.Cells.Item($row, 2).Value = "21:32:59"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "03"

foreach ($col in $cols) {
    $ws.Cells.Item($row, $col).ClearFormats()
}

$ws.Cells.Item($row, 5).Value = 134625
$ws.Cells.Item($row, 6).Value = 141690
$ws.Cells.Item($row, 7).Value = 171610
$ws.Cells.Item($row, 8).Value = 149214
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 123299
$ws.Cells.Item($row, 11).Value = 223969
$ws.Cells.Item($row, 12).Value = 256736
$ws.Cells.Item($row, 13).Value = 185286
$ws.Cells.Item($row, 14).Value = 110036
$ws.Cells.Item($row, 15).Value = 41332
$ws.Cells.Item($row, 16).Value = 30891
$ws.Cells.Item($row, 17).Value = 73538
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42483
$ws.Cells.Item($row, 20).Value = -1
